$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns for the "ward" bed metrics, between the existing
# HDU columns (B-G) and the cumulative-deaths columns (old H, I).
# This automatically shifts the old H/I ("Cumulative deaths", "Deaths SD")
# data into K/L, keeping their values/styles/shared-strings intact.
$ws.Range("H1:J1").EntireColumn.Insert()

# New header labels for the inserted "ward" columns.
$ws.Range("H1").Value2 = "Peak ward bed needs timing"
$ws.Range("I1").Value2 = "Mean peak ward bed needs"
$ws.Range("J1").Value2 = "Ward SD"

# Row 2 - Base (0%)
$ws.Range("A2").Value2 = "Base (0%)"
$ws.Range("B2").Value = [DateTime]"2020-05-07"
$ws.Range("C2").Value2 = 1807
$ws.Range("D2").Value2 = 30
$ws.Range("E2").Value = [DateTime]"2020-05-13"
$ws.Range("F2").Value2 = 978
$ws.Range("G2").Value2 = 26
$ws.Range("B2").Copy($ws.Range("H2"))
$ws.Range("H2").Value = [DateTime]"2020-05-25"
$ws.Range("I2").Value2 = 1439
$ws.Range("J2").Value2 = 33
$ws.Range("K2").Value2 = 3012
$ws.Range("L2").Value2 = 473

# Row 3 - 20% reduction
$ws.Range("A3").Value2 = "20% reduction"
$ws.Range("B3").Value = [DateTime]"2020-05-09"
$ws.Range("C3").Value2 = 1766
$ws.Range("D3").Value2 = 29
$ws.Range("E3").Value = [DateTime]"2020-05-15"
$ws.Range("F3").Value2 = 959
$ws.Range("G3").Value2 = 27
$ws.Range("B3").Copy($ws.Range("H3"))
$ws.Range("H3").Value = [DateTime]"2020-05-26"
$ws.Range("I3").Value2 = 1417
$ws.Range("J3").Value2 = 32
$ws.Range("K3").Value2 = 3011
$ws.Range("L3").Value2 = 476

# Row 4 - 40% reduction
$ws.Range("A4").Value2 = "40% reduction"
$ws.Range("B4").Value = [DateTime]"2020-05-09"
$ws.Range("C4").Value2 = 1773
$ws.Range("D4").Value2 = 28
$ws.Range("E4").Value = [DateTime]"2020-05-14"
$ws.Range("F4").Value2 = 962
$ws.Range("G4").Value2 = 29
$ws.Range("B4").Copy($ws.Range("H4"))
$ws.Range("H4").Value = [DateTime]"2020-05-26"
$ws.Range("I4").Value2 = 1420
$ws.Range("J4").Value2 = 33
$ws.Range("K4").Value2 = 3010
$ws.Range("L4").Value2 = 478

# Row 5 - 50% reduction
$ws.Range("A5").Value2 = "50% reduction"
$ws.Range("B5").Value = [DateTime]"2020-05-17"
$ws.Range("C5").Value2 = 1548
$ws.Range("D5").Value2 = 28
$ws.Range("E5").Value = [DateTime]"2020-05-22"
$ws.Range("F5").Value2 = 868
$ws.Range("G5").Value2 = 27
$ws.Range("B5").Copy($ws.Range("H5"))
$ws.Range("H5").Value = [DateTime]"2020-06-03"
$ws.Range("I5").Value2 = 1291
$ws.Range("J5").Value2 = 32
$ws.Range("K5").Value2 = 2992
$ws.Range("L5").Value2 = 511

# Row 6 - 60% reduction
$ws.Range("A6").Value2 = "60% reduction"
$ws.Range("B6").Value = [DateTime]"2020-05-24"
$ws.Range("C6").Value2 = 1375
$ws.Range("D6").Value2 = 27
$ws.Range("E6").Value = [DateTime]"2020-05-30"
$ws.Range("F6").Value2 = 788
$ws.Range("G6").Value2 = 25
$ws.Range("B6").Copy($ws.Range("H6"))
$ws.Range("H6").Value = [DateTime]"2020-06-11"
$ws.Range("I6").Value2 = 1181
$ws.Range("J6").Value2 = 32
$ws.Range("K6").Value2 = 2959
$ws.Range("L6").Value2 = 531

# Row 7 - 73% reduction
$ws.Range("A7").Value2 = "73% reduction"
$ws.Range("B7").Value = [DateTime]"2020-06-08"
$ws.Range("C7").Value2 = 1047
$ws.Range("D7").Value2 = 22
$ws.Range("E7").Value = [DateTime]"2020-06-13"
$ws.Range("F7").Value2 = 621
$ws.Range("G7").Value2 = 22
$ws.Range("B7").Copy($ws.Range("H7"))
$ws.Range("H7").Value = [DateTime]"2020-06-25"
$ws.Range("I7").Value2 = 947
$ws.Range("J7").Value2 = 28
$ws.Range("K7").Value2 = 2880
$ws.Range("L7").Value2 = 602
